# "took test resume out of protected view"
#
# When Word takes a document out of Protected View it re-saves it, and as
# part of that the internal "_GoBack" bookmark (which always marks the
# location of the most recent edit) is re-anchored to the very start of
# the document's main story. Recreate that here: move the existing
# "_GoBack" bookmark from wherever it currently sits onto a
# zero-length range at the very start of the first paragraph.

$d = $word.ActiveDocument

# The engine's Bookmarks.Add has trouble inserting a *collapsed* bookmark
# exactly at absolute document position 0 when that position also happens
# to be the start of the very first paragraph in the body. Work around it
# by briefly inserting a throwaway paragraph in front of everything so the
# real first paragraph is (temporarily) paragraph #2, placing the bookmark
# at that paragraph's start (which behaves correctly), and then deleting
# the scratch paragraph again. Bookmarks survive the deletion of the
# paragraph that precedes them and keep their resolved position.

$frontRange = $d.Range(0, 0)
$frontRange.InsertParagraphBefore()

$targetStart = $d.Paragraphs(2).Range.Start
$targetRange = $d.Range($targetStart, $targetStart)

# Adding a bookmark named "_GoBack" automatically replaces any existing
# bookmark of that name elsewhere in the document (bookmark names are
# unique), so this both relocates it and removes the old occurrence in a
# single step.
$d.Bookmarks.Add("_GoBack", $targetRange)

# Remove the scratch paragraph used purely to dodge the position-0 quirk.
$d.Paragraphs(1).Range.Delete()
